$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Stat Query" text in C2 and C3: add the node alias "a" to the
# arm pattern so downstream steps can reference the arm (a:arm) instead of
# the anonymous (:arm).
$newStatQuery = "MATCH (s:specimen)-->(c:case)-->(a:arm)-->(ct:clinical_trial)`r`n    WHERE c.gender = `"MALE`"`r`nOPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)`r`nRETURN `r`n     COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,`r`n     COUNT(DISTINCT c.case_id) AS Cases,`r`n      COUNT(DISTINCT f) AS Files"

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery

# Move the active selection from D3 to B3.
$ws.Range("A3").Activate()
$ws.Range("B3").Select()
